$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column M header cell (style-only, matches J1:L1 header style), row 1 ---
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 3: add "Input or Output" label in A3 (Frequency stays in B3) ---
$ws.Range("A3").Value = "Input or Output"

# --- Row 4: add "Input or Output" label in A4, and push B4's old "Heat" value ---
# down into a brand-new row 5 renamed to "Sound", while B4 itself becomes the
# brand-new "Temperature" value.
$ws.Range("A4").Value = "Input or Output"
$ws.Range("A5").Value = "Input or Output"
$ws.Range("B5").Value = "Heat"
$ws.Range("B4").Value = "Temperature"
$ws.Range("B5").Value = "Sound"

# --- Update the active selection to match the authored state ---
$ws.Range("L4").Select()
